$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 and 4 (2019 B/C) are swapped
$ws.Range("A3").Value = "2019年C"
$ws.Range("B3").Value = 3.8
$ws.Range("C3").Value = 1.4
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = 1.5
$ws.Range("F3").Value = 7

$ws.Range("A4").Value = "2019年B"
$ws.Range("B4").Value = 5.3
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = 1.5
$ws.Range("F4").Value = 7.4

# Rows 7 and 8 (2020 B/C) are swapped
$ws.Range("A7").Value = "2020年C"
$ws.Range("B7").Value = -9.199999999999999
$ws.Range("C7").Value = -6.8
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = -16.6
$ws.Range("F7").Value = 2.9

$ws.Range("A8").Value = "2020年B"
$ws.Range("B8").Value = -14.1
$ws.Range("C8").Value = -14.4
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = -23.6
$ws.Range("F8").Value = 3.8

# Rows 11 and 12 (2021 B/C) are swapped (note D column gains/loses a value)
$ws.Range("A11").Value = "2021年C"
$ws.Range("B11").Value = 17.1
$ws.Range("C11").Value = 19.8
$ws.Range("D11").Value = 14.5
$ws.Range("E11").Value = 19.9
$ws.Range("F11").Value = 6.2

$ws.Range("A12").Value = "2021年B"
$ws.Range("B12").Value = 24.6
$ws.Range("C12").Value = 32.4
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = 29.5
$ws.Range("F12").Value = 6.2
